# Comprobación de que la casilla de nombre se obtiene correctamente
#
# Adds a header row above the existing data table, labelling each of the
# 18 columns, and formats that header row (bold, centred, thin box border
# around every cell). All pre-existing data rows (and their hyperlinks on
# column F) shift down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$urls = @(
    "http://evidentia.test/20/profiles/view/1",
    "http://evidentia.test/20/profiles/view/2",
    "http://evidentia.test/20/profiles/view/3",
    "http://evidentia.test/20/profiles/view/4",
    "http://evidentia.test/20/profiles/view/5",
    "http://evidentia.test/20/profiles/view/6",
    "http://evidentia.test/20/profiles/view/7"
)

# The hyperlink collection is worksheet-wide; clear it before shifting rows
# so we don't leave stale refs behind, then rebuild it afterwards.
$ws.Range("A1:R7").Hyperlinks.Delete()

# Push the existing data down and make room for the new header row.
$ws.Rows.Item(1).Insert()

$headers = @(
    "DNI",
    "Apellidos",
    "Nombre",
    "Uvus",
    "Correo",
    "Perfil",
    "Participación",
    "Comité",
    "Evidencia aleatoria",
    "Horas de evidencia aleatoria",
    "Eventos asistidos",
    "Horas de asistencia",
    "Reuniones asistidas",
    "Bono de horas",
    "Horas de reuniones",
    "Evidencias registradas",
    "Horas de evidencias",
    "Horas en total"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Re-create the hyperlinks on column F, now at rows 2..8.
for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = $i + 2
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 6), $urls[$i])
}

# Build the header look (bold, thin border all round, centred/top aligned)
# on a scratch cell once, then format-paint it onto the header row. Doing
# it this way (one template build + a single paste) keeps the style table
# tidy instead of growing a new style per cell.
$template = $ws.Cells.Item(100, 26)
$template.Font.Bold = $true
$template.Borders.LineStyle = 1
$template.HorizontalAlignment = -4108
$template.VerticalAlignment = -4160

$template.Copy()
$ws.Range("A1:R1").PasteSpecial(-4122)
$template.Clear()
$excel.CutCopyMode = $false
